{"js": "// Replace the date header and every \"A\u00d7B=C\" answer cell with the new\n// values from the latest run. Order matters: \"733\u00d79=6597\" is both an old\n// value (row 5 table, col 4) and a new value (row 1 table, col 4) in this\n// pass, so the row-5 replacement is performed first to avoid temporarily\n// creating a duplicate / mismatched match.\nconst replacements = [\n  [\"733\u00d79=6597\", \"561\u00d78=4488\"],\n  [\"2025-03-29 Saturday\", \"2025-03-30 Sunday\"],\n  [\"896\u00d72=1792\", \"616\u00d77=4312\"],\n  [\"620\u00d72=1240\", \"909\u00d78=7272\"],\n  [\"992\u00d72=1984\", \"215\u00d78=1720\"],\n  [\"684\u00d77=4788\", \"733\u00d79=6597\"],\n  [\"219\u00d79=1971\", \"875\u00d77=6125\"],\n  [\"128\u00d72=256\", \"747\u00d78=5976\"],\n  [\"473\u00d77=3311\", \"334\u00d79=3006\"],\n  [\"592\u00d74=2368\", \"695\u00d79=6255\"],\n  [\"492\u00d72=984\", \"863\u00d72=1726\"],\n  [\"881\u00d79=7929\", \"195\u00d78=1560\"],\n  [\"108\u00d74=432\", \"643\u00d74=2572\"],\n  [\"987\u00d74=3948\", \"129\u00d77=903\"],\n  [\"494\u00d77=3458\", \"387\u00d78=3096\"],\n  [\"279\u00d77=1953\", \"449\u00d73=1347\"],\n  [\"595\u00d74=2380\", \"146\u00d73=438\"],\n  [\"391\u00d74=1564\", \"319\u00d76=1914\"],\n  [\"442\u00d74=1768\", \"348\u00d76=2088\"],\n  [\"170\u00d79=1530\", \"584\u00d78=4672\"],\n  [\"945\u00d78=7560\", \"122\u00d74=488\"],\n  [\"460\u00d74=1840\", \"400\u00d72=800\"],\n  [\"297\u00d77=2079\", \"265\u00d78=2120\"],\n  [\"395\u00d74=1580\", \"959\u00d75=4795\"],\n  [\"159\u00d73=477\", \"339\u00d73=1017\"],\n  [\"982\u00d74=3928\", \"410\u00d75=2050\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date header and every \"A\u00d7B=C\" answer cell with the new\n# values from the latest run. Order matters: \"733\u00d79=6597\" is both an old\n# value (row 5 of the table) and a new value (row 1 of the table) in this\n# pass, so the row-5 replacement is performed first to avoid temporarily\n# creating an ambiguous duplicate match.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"733\u00d79=6597\", \"561\u00d78=4488\"),\n  @(\"2025-03-29 Saturday\", \"2025-03-30 Sunday\"),\n  @(\"896\u00d72=1792\", \"616\u00d77=4312\"),\n  @(\"620\u00d72=1240\", \"909\u00d78=7272\"),\n  @(\"992\u00d72=1984\", \"215\u00d78=1720\"),\n  @(\"684\u00d77=4788\", \"733\u00d79=6597\"),\n  @(\"219\u00d79=1971\", \"875\u00d77=6125\"),\n  @(\"128\u00d72=256\", \"747\u00d78=5976\"),\n  @(\"473\u00d77=3311\", \"334\u00d79=3006\"),\n  @(\"592\u00d74=2368\", \"695\u00d79=6255\"),\n  @(\"492\u00d72=984\", \"863\u00d72=1726\"),\n  @(\"881\u00d79=7929\", \"195\u00d78=1560\"),\n  @(\"108\u00d74=432\", \"643\u00d74=2572\"),\n  @(\"987\u00d74=3948\", \"129\u00d77=903\"),\n  @(\"494\u00d77=3458\", \"387\u00d78=3096\"),\n  @(\"279\u00d77=1953\", \"449\u00d73=1347\"),\n  @(\"595\u00d74=2380\", \"146\u00d73=438\"),\n  @(\"391\u00d74=1564\", \"319\u00d76=1914\"),\n  @(\"442\u00d74=1768\", \"348\u00d76=2088\"),\n  @(\"170\u00d79=1530\", \"584\u00d78=4672\"),\n  @(\"945\u00d78=7560\", \"122\u00d74=488\"),\n  @(\"460\u00d74=1840\", \"400\u00d72=800\"),\n  @(\"297\u00d77=2079\", \"265\u00d78=2120\"),\n  @(\"395\u00d74=1580\", \"959\u00d75=4795\"),\n  @(\"159\u00d73=477\", \"339\u00d73=1017\"),\n  @(\"982\u00d74=3928\", \"410\u00d75=2050\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute(\n    $oldText,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $newText,\n    2\n  ) | Out-Null\n}\n"}
